$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed coin-ranking data.
# Cells are stored as plain text, so force a text number-format before assigning
# to avoid Excel's automatic number/percentage auto-conversion mangling the values
# (e.g. dropping significant trailing zeros like '0.001507' -> 0.001507 as a float).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "13.53%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05716"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.99%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.558"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.30%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8575"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8729"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.88%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1340"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.68%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06888"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.71%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02861"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.24%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09391"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.13%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001507"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.43%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04171"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-9.28%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006008"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-93.97%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005992"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.89%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.513"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.87%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.012"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.40%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.177"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.35%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3146"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.11%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03377"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.37%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1302"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.603"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.83%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1374"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001209"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.50%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004478"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.30%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001178"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.63%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001386"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.96%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03768"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.53%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005778"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.53%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.40%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002274"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009673"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "19.37%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005074"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.22%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.25%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08967"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.97%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002757"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "12.17%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002095"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.25%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001996"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.25%"
